$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Language")
$ws.Range("B2:B52").NumberFormat = "@"

$ws.Cells.Item(2, 1).Value = "English"
$ws.Cells.Item(2, 2).Value = "1203"
$ws.Cells.Item(3, 1).Value = "Spanish"
$ws.Cells.Item(3, 2).Value = "488"
$ws.Cells.Item(4, 1).Value = "Russian"
$ws.Cells.Item(4, 2).Value = "421"
$ws.Cells.Item(5, 1).Value = "French"
$ws.Cells.Item(5, 2).Value = "390"
$ws.Cells.Item(6, 1).Value = "Portuguese"
$ws.Cells.Item(6, 2).Value = "312"
$ws.Cells.Item(7, 1).Value = "Arabic"
$ws.Cells.Item(7, 2).Value = "171"
$ws.Cells.Item(8, 1).Value = "Chinese (China)"
$ws.Cells.Item(8, 2).Value = "137"
$ws.Cells.Item(9, 1).Value = "Korean"
$ws.Cells.Item(9, 2).Value = "137"
$ws.Cells.Item(10, 1).Value = "German"
$ws.Cells.Item(10, 2).Value = "110"
$ws.Cells.Item(11, 1).Value = "Vietnamese"
$ws.Cells.Item(11, 2).Value = "110"
$ws.Cells.Item(12, 1).Value = "Portuguese (Brazil)"
$ws.Cells.Item(12, 2).Value = "73"
$ws.Cells.Item(13, 1).Value = "Turkish"
$ws.Cells.Item(13, 2).Value = "70"
$ws.Cells.Item(14, 1).Value = "Japanese"
$ws.Cells.Item(14, 2).Value = "67"
$ws.Cells.Item(15, 1).Value = "Italian"
$ws.Cells.Item(15, 2).Value = "45"
$ws.Cells.Item(16, 1).Value = "Persian"
$ws.Cells.Item(16, 2).Value = "35"
$ws.Cells.Item(17, 1).Value = "Chinese (Traditional)"
$ws.Cells.Item(17, 2).Value = "31"
$ws.Cells.Item(18, 1).Value = "Chinese"
$ws.Cells.Item(18, 2).Value = "30"
$ws.Cells.Item(19, 1).Value = "Greek"
$ws.Cells.Item(19, 2).Value = "25"
$ws.Cells.Item(20, 1).Value = "Ukrainian"
$ws.Cells.Item(20, 2).Value = "22"
$ws.Cells.Item(21, 1).Value = "Hindi"
$ws.Cells.Item(21, 2).Value = "20"
$ws.Cells.Item(22, 1).Value = "Romanian"
$ws.Cells.Item(22, 2).Value = "16"
$ws.Cells.Item(23, 1).Value = "Hebrew"
$ws.Cells.Item(23, 2).Value = "13"
$ws.Cells.Item(24, 1).Value = "Catalan"
$ws.Cells.Item(24, 2).Value = "11"
$ws.Cells.Item(25, 1).Value = "Dutch"
$ws.Cells.Item(25, 2).Value = "10"
$ws.Cells.Item(26, 1).Value = "Thai"
$ws.Cells.Item(26, 2).Value = "9"
$ws.Cells.Item(27, 1).Value = "Polish"
$ws.Cells.Item(27, 2).Value = "8"
$ws.Cells.Item(28, 1).Value = "Afrikaans"
$ws.Cells.Item(28, 2).Value = "7"
$ws.Cells.Item(29, 1).Value = "Mongolian"
$ws.Cells.Item(29, 2).Value = "7"
$ws.Cells.Item(30, 1).Value = "Serbian"
$ws.Cells.Item(30, 2).Value = "7"
$ws.Cells.Item(31, 1).Value = "Indonesian"
$ws.Cells.Item(31, 2).Value = "6"
$ws.Cells.Item(32, 1).Value = "Portuguese (Portugal)"
$ws.Cells.Item(32, 2).Value = "6"
$ws.Cells.Item(33, 1).Value = "Tamil"
$ws.Cells.Item(33, 2).Value = "6"
$ws.Cells.Item(34, 1).Value = "Hungarian"
$ws.Cells.Item(34, 2).Value = "5"
$ws.Cells.Item(35, 1).Value = "Slovak"
$ws.Cells.Item(35, 2).Value = "5"
$ws.Cells.Item(36, 1).Value = "Bengali"
$ws.Cells.Item(36, 2).Value = "4"
$ws.Cells.Item(37, 1).Value = "Czech"
$ws.Cells.Item(37, 2).Value = "4"
$ws.Cells.Item(38, 1).Value = "Telugu"
$ws.Cells.Item(38, 2).Value = "4"
$ws.Cells.Item(39, 1).Value = "Urdu"
$ws.Cells.Item(39, 2).Value = "3"
$ws.Cells.Item(40, 1).Value = "Albanian"
$ws.Cells.Item(40, 2).Value = "2"
$ws.Cells.Item(41, 1).Value = "Bulgarian"
$ws.Cells.Item(41, 2).Value = "2"
$ws.Cells.Item(42, 1).Value = "Burmese"
$ws.Cells.Item(42, 2).Value = "2"
$ws.Cells.Item(43, 1).Value = "Estonian"
$ws.Cells.Item(43, 2).Value = "2"
$ws.Cells.Item(44, 1).Value = "Georgian"
$ws.Cells.Item(44, 2).Value = "2"
$ws.Cells.Item(45, 1).Value = "Kazakh"
$ws.Cells.Item(45, 2).Value = "2"
$ws.Cells.Item(46, 1).Value = "Lithuanian"
$ws.Cells.Item(46, 2).Value = "2"
$ws.Cells.Item(47, 1).Value = "Uzbek"
$ws.Cells.Item(47, 2).Value = "2"
$ws.Cells.Item(48, 1).Value = "Croatian"
$ws.Cells.Item(48, 2).Value = "1"
$ws.Cells.Item(49, 1).Value = "Javanese"
$ws.Cells.Item(49, 2).Value = "1"
$ws.Cells.Item(50, 1).Value = "Marathi"
$ws.Cells.Item(50, 2).Value = "1"
$ws.Cells.Item(51, 1).Value = "Swedish"
$ws.Cells.Item(51, 2).Value = "1"
$ws.Cells.Item(52, 1).Value = "TOTAL COURSES"
$ws.Cells.Item(52, 2).Value = "4048"

$ws3 = $wb.Worksheets.Item("Level")
$ws3.Range("B2:B6").NumberFormat = "@"

$ws3.Cells.Item(2, 1).Value = "Intermediate"
$ws3.Cells.Item(2, 2).Value = "530"
$ws3.Cells.Item(3, 1).Value = "Beginner"
$ws3.Cells.Item(3, 2).Value = "492"
$ws3.Cells.Item(4, 1).Value = "Mixed"
$ws3.Cells.Item(4, 2).Value = "230"
$ws3.Cells.Item(5, 1).Value = "Advanced"
$ws3.Cells.Item(5, 2).Value = "78"
$ws3.Cells.Item(6, 1).Value = "TOTAL COURSES"
$ws3.Cells.Item(6, 2).Value = "1330"
